$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$val) {
    # Cells in this sheet store prices/percentages as TEXT (t="inlineStr").
    # Excel auto-converts plain-decimal-looking strings to numbers on .Value
    # assignment, so force those through the text quote-prefix, then strip the
    # resulting "Normal-ish" style back down so no stray number-format sticks.
    if ($val -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $val
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}

Set-TextValue $ws.Range("D2") "69.666.43"
Set-TextValue $ws.Range("E2") "  +2.59%  "
Set-TextValue $ws.Range("D3") "3.378.85"
Set-TextValue $ws.Range("E3") "  +3.60%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "191.11"
Set-TextValue $ws.Range("E5") "  +2.94%  "
Set-TextValue $ws.Range("D6") "593.22"
Set-TextValue $ws.Range("E6") "  +1.96%  "
Set-TextValue $ws.Range("E7") "  +0.04%  "
Set-TextValue $ws.Range("E8") "  +0.82%  "
Set-TextValue $ws.Range("E9") "  +1.93%  "
Set-TextValue $ws.Range("E10") "  +2.63%  "
Set-TextValue $ws.Range("E11") "  +1.66%  "
Set-TextValue $ws.Range("D12") "3.970.11"
Set-TextValue $ws.Range("E12") "  +3.80%  "
Set-TextValue $ws.Range("E13") "  -0.63%  "
Set-TextValue $ws.Range("D15") "69.657.09"
Set-TextValue $ws.Range("E15") "  +2.70%  "
Set-TextValue $ws.Range("D17") "3.382.98"
Set-TextValue $ws.Range("E17") "  +3.01%  "
Set-TextValue $ws.Range("D18") "453.60"
Set-TextValue $ws.Range("E18") "  +15.14%  "
Set-TextValue $ws.Range("E19") "  +1.26%  "
Set-TextValue $ws.Range("D20") "13.83"
Set-TextValue $ws.Range("E20") "  +1.68%  "
Set-TextValue $ws.Range("D21") "7.78"
Set-TextValue $ws.Range("E21") "  +2.03%  "
Set-TextValue $ws.Range("D22") "75.91"
Set-TextValue $ws.Range("E22") "  +5.97%  "
Set-TextValue $ws.Range("E23") "  +0.06%  "
Set-TextValue $ws.Range("D24") "0.521"
Set-TextValue $ws.Range("E24") "  +0.96%  "
Set-TextValue $ws.Range("E25") "  +3.03%  "
Set-TextValue $ws.Range("E26") "  +1.84%  "
Set-TextValue $ws.Range("E27") "  -0.72%  "
Set-TextValue $ws.Range("E28") "  -0.14%  "
Set-TextValue $ws.Range("E29") "  +3.18%  "
Set-TextValue $ws.Range("D30") "23.39"
Set-TextValue $ws.Range("E30") "  +2.99%  "
Set-TextValue $ws.Range("E31") "  +0.99%  "
Set-TextValue $ws.Range("E32") "  +2.31%  "
Set-TextValue $ws.Range("E33") "  +0.12%  "
Set-TextValue $ws.Range("E34") "  -0.05%  "
Set-TextValue $ws.Range("E35") "  +6.26%  "
Set-TextValue $ws.Range("D36") "164.57"
Set-TextValue $ws.Range("E36") "  +1.22%  "
Set-TextValue $ws.Range("E37") "  +2.21%  "
Set-TextValue $ws.Range("D38") "27.76"
Set-TextValue $ws.Range("E38") "  +3.41%  "
Set-TextValue $ws.Range("E39") "  +0.35%  "
Set-TextValue $ws.Range("D40") "4.59"
Set-TextValue $ws.Range("E40") "  +1.03%  "
Set-TextValue $ws.Range("D41") "6.59"
Set-TextValue $ws.Range("E41") "  +1.67%  "
Set-TextValue $ws.Range("D42") "2.744.02"
Set-TextValue $ws.Range("E43") "  +1.93%  "
Set-TextValue $ws.Range("E44") "  +2.67%  "
Set-TextValue $ws.Range("B45") "Hedera"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D45") "0.0689"
Set-TextValue $ws.Range("E45") "  -0.07%  "
Set-TextValue $ws.Range("B46") "OKB"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D46") "41.10"
Set-TextValue $ws.Range("E46") "  +1.09%  "
Set-TextValue $ws.Range("D47") "339.33"
Set-TextValue $ws.Range("E47") "  +1.57%  "
Set-TextValue $ws.Range("D48") "0.0284"
Set-TextValue $ws.Range("E48") "  +2.47%  "
Set-TextValue $ws.Range("D49") "32.86"
Set-TextValue $ws.Range("E49") "  +6.33%  "
Set-TextValue $ws.Range("E50") "  +4.49%  "
Set-TextValue $ws.Range("E51") "  -0.58%  "
